{"js": "// Remove the site-footer paragraphs that trail the \"Requisitos\" section:\n//   - the blank spacer paragraph right after \"LOQ4205: Sistemas Produtivos II (Requisito fraco)\"\n//   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n// while leaving the following blank paragraph (and the page-break paragraph after it) untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nconst requisitoIdx = items.findIndex(p =>\n    p.text.indexOf(\"LOQ4205: Sistemas Produtivos II (Requisito fraco)\") !== -1\n);\nif (requisitoIdx === -1) {\n    throw new Error(\"Could not locate the 'LOQ4205: ... (Requisito fraco)' paragraph\");\n}\n\nconst jupiterIdx = items.findIndex(p =>\n    p.text.indexOf(\"Ver no Jupiter Salvar em pdf Salvar em docx\") !== -1\n);\nif (jupiterIdx === -1) {\n    throw new Error(\"Could not locate the 'Ver no Jupiter...' paragraph\");\n}\n\nconst copyrightIdx = items.findIndex(p => p.text.indexOf(\"Contact: luizeleno@usp.br\") !== -1);\nif (copyrightIdx === -1) {\n    throw new Error(\"Could not locate the copyright/footer paragraph\");\n}\n\n// The blank spacer paragraph sits immediately before the \"Ver no Jupiter...\" paragraph.\nconst blankIdx = jupiterIdx - 1;\nif (blankIdx <= requisitoIdx || items[blankIdx].text !== \"\") {\n    throw new Error(\"Unexpected document structure around the footer paragraphs\");\n}\n\n// Delete from bottom to top so earlier indices stay valid.\nitems[copyrightIdx].delete();\nitems[jupiterIdx].delete();\nitems[blankIdx].delete();\n\nawait context.sync();\n", "ps1": "# Remove the site-footer paragraphs that trail the \"Requisitos\" section:\n#   - the blank spacer paragraph right after \"LOQ4205: Sistemas Produtivos II (Requisito fraco)\"\n#   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   - \"... Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n# while leaving the following blank paragraph (and the page-break paragraph after it) untouched.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndexContaining($doc, $needle) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        if ($doc.Paragraphs.Item($i).Range.Text -like \"*$needle*\") {\n            return $i\n        }\n    }\n    return -1\n}\n\nfunction Remove-ParagraphContaining($doc, $needle) {\n    $rng = $doc.Content\n    $found = $rng.Find.Execute($needle)\n    if (-not $found) {\n        throw \"Could not find paragraph containing: $needle\"\n    }\n    [void]$rng.Expand(4)  # wdParagraph\n    [void]$rng.Delete()\n}\n\n# Delete the two non-blank footer paragraphs by locating their text.\nRemove-ParagraphContaining $d \"Contact: luizeleno@usp.br\"\nRemove-ParagraphContaining $d \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n\n# The blank spacer paragraph has no text to search for; it directly follows\n# the \"Requisitos\" line, so locate it relative to that paragraph.\n$reqIdx = Find-ParagraphIndexContaining $d \"LOQ4205: Sistemas Produtivos II (Requisito fraco)\"\nif ($reqIdx -eq -1) {\n    throw \"Could not locate the 'LOQ4205: ... (Requisito fraco)' paragraph\"\n}\n\n$blankIdx = $reqIdx + 1\n$blankPara = $d.Paragraphs.Item($blankIdx)\nif ($blankPara.Range.Text.Trim() -ne \"\") {\n    throw \"Unexpected document structure around the footer paragraphs\"\n}\n[void]$blankPara.Range.Delete()\n"}
